$d = $word.ActiveDocument

# 1) The "Autres : ..." paragraph must disappear and the following "Langages : ..."
#    paragraph must become "Langages : python, matlab, c, c++" -- i.e. the two
#    paragraphs collapse into one. The "Autres" paragraph sits immediately after a
#    section-break-only (empty) paragraph, so we must not delete that paragraph's
#    own mark directly (doing so merges/loses the section break). Instead we
#    overwrite the "Autres" paragraph's text in place with the final "Langages"
#    text, then remove the old, now-duplicate "Langages" paragraph that follows it
#    (which is not adjacent to the section break, so it is safe to delete outright).
$d.Content.Find.Execute(
    "Autres : marketing, google analytics, internes comme externes, presse, affichage, site centric, formats",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Langages : python, matlab, c, c++", 2) | Out-Null

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Langages : r, python, matlab, c, c++*") {
        [void]$p.Range.Delete()
        break
    }
}

# 2) Update the remaining skill-category lines' text in place.
$d.Content.Find.Execute(
    "Visualisation : web analytics, tableau", $true, $false, $false, $false, $false,
    $true, 1, $false, "Data Science : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn", 2) | Out-Null

$d.Content.Find.Execute(
    "MLOps : Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit", $true, $false, $false, $false, $false,
    $true, 1, $false, "Visualisation : tableau", 2) | Out-Null

$d.Content.Find.Execute(
    "ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn", $true, $false, $false, $false, $false,
    $true, 1, $false, "Machine Learning : Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit", 2) | Out-Null
